$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-03 Wednesday" "2024-01-04 Thursday"

Replace-Text "21×36=" "78×60="
Replace-Text "27×50=" "30×84="
Replace-Text "90×29=" "84×55="
Replace-Text "30×57=" "96×39="
Replace-Text "97×49=" "38×39="
Replace-Text "69×29=" "86×95="
Replace-Text "54×46=" "28×30="
Replace-Text "40×23=" "67×29="
Replace-Text "32×37=" "91×32="
Replace-Text "89×15=" "72×72="
Replace-Text "55×59=" "67×13="
Replace-Text "77×88=" "66×94="
Replace-Text "82×30=" "47×84="
Replace-Text "69×14=" "31×22="
Replace-Text "30×38=" "44×20="
Replace-Text "60×91=" "81×97="
Replace-Text "33×28=" "98×36="
Replace-Text "16×58=" "81×88="
Replace-Text "27×26=" "59×36="
Replace-Text "32×68=" "75×87="
Replace-Text "92×55=" "78×46="
Replace-Text "72×50=" "24×75="
Replace-Text "77×14=" "57×30="
Replace-Text "67×64=" "16×27="
Replace-Text "57×72=" "51×95="
